$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 14:42"

# Update country stats that changed between snapshots, then the table
# (sorted descending by "Casos totales") is re-sorted, which is why some
# rows below show a different country than before even though that
# country's own figures did not change.

# Row 6: 'Estados Unidos' -> 'Estados Unidos'
$ws.Range("B6").Value = 68814
$ws.Range("C6").Value = 603
$ws.Range("E6").Value = 67349

# Row 8: 'Alemania' -> 'Alemania'
$ws.Range("B8").Value = 39572
$ws.Range("C8").Value = 2249
$ws.Range("D8").Value = 3959
$ws.Range("E8").Value = 35389
$ws.Range("G8").Value = 18
$ws.Range("H8").Value = 224

# Row 19: 'Noruega' -> 'Noruega'
$ws.Range("B19").Value = 3250
$ws.Range("C19").Value = 166
$ws.Range("E19").Value = 3230

# Row 20: 'Australia' -> 'Suecia'
$ws.Range("A20").Value = "Suecia"
$ws.Range("B20").Value = 2840
$ws.Range("C20").Value = 314
$ws.Range("D20").Value = 16
$ws.Range("E20").Value = 2758
$ws.Range("F20").Value = 176
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 66

# Row 21: 'Israel' -> 'Australia'
$ws.Range("A21").Value = "Australia"
$ws.Range("B21").Value = 2806
$ws.Range("C21").Value = 130
$ws.Range("D21").Value = 170
$ws.Range("E21").Value = 2623
$ws.Range("F21").Value = 11
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = 13

# Row 22: 'Brasil' -> 'Israel'
$ws.Range("A22").Value = "Israel"
$ws.Range("B22").Value = 2666
$ws.Range("C22").Value = 297
$ws.Range("D22").Value = 68
$ws.Range("E22").Value = 2590
$ws.Range("F22").Value = 39
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 8

# Row 23: 'Suecia' -> 'Brasil'
$ws.Range("A23").Value = "Brasil"
$ws.Range("B23").Value = 2567
$ws.Range("C23").Value = 13
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = 2500
$ws.Range("F23").Value = 18
$ws.Range("H23").Value = 61

# Row 25: 'Malasia' -> 'Malasia'
$ws.Range("D25").Value = 215
$ws.Range("E25").Value = 1793

# Row 31: 'Ecuador' -> 'Chile'
$ws.Range("A31").Value = "Chile"
$ws.Range("B31").Value = 1306
$ws.Range("C31").Value = 164
$ws.Range("D31").Value = 22
$ws.Range("E31").Value = 1280
$ws.Range("F31").Value = 7
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 4

# Row 32: 'Chile' -> 'Ecuador'
$ws.Range("A32").Value = "Ecuador"
$ws.Range("B32").Value = 1211
$ws.Range("D32").Value = 3
$ws.Range("E32").Value = 1179
$ws.Range("F32").Value = 2
$ws.Range("H32").Value = 29

# Row 33: 'Pakistan' -> 'Pakistan'
$ws.Range("B33").Value = 1128
$ws.Range("C33").Value = 65
$ws.Range("E33").Value = 1099

# Row 38: 'Finlandia' -> 'Finlandia'
$ws.Range("F38").Value = 24

# Row 47: 'Singapur' -> 'Singapur'
$ws.Range("B47").Value = 683
$ws.Range("C47").Value = 52
$ws.Range("D47").Value = 172
$ws.Range("E47").Value = 509

# Row 100: 'Sri Lanka' -> 'Sri Lanka'
$ws.Range("B100").Value = 104
$ws.Range("C100").Value = 2
$ws.Range("D100").Value = 6
$ws.Range("E100").Value = 98
